$wb = $excel.ActiveWorkbook

# Sheet-specific "container name" helper text (H3) and example container name
$helpText = @{
    "Tank"   = "Optional. Tank name of the container this group should be in. If there is already an equivalent group in that tank, they will be merged. E.g. LP1";
    "Trough" = "Optional. Trough name of the container this group should be in. If there is already an equivalent group in that trough, they will be merged. E.g. TR1";
    "Drawer" = "Optional. Drawer name of the container this group should be in. If there is already an equivalent group in that Drawer, they will be merged. E.g. HU1.2";
    "Cup"    = "Optional. Cup name of the container this group should be in. If there is already an equivalent group in that cup, they will be merged. E.g. HU1.2.3";
}

$sheetOrder = @("Tank", "Trough", "Drawer", "Cup")

foreach ($name in $sheetOrder) {
    $ws = $wb.Worksheets.Item($name)

    # update the parser-specific container-name help text
    $ws.Range("H3").Value = $helpText[$name]

    # move the tab/cell selection onto the new important column
    $ws.Activate()
    $ws.Range("H4").Select()
}

# the Cup instructions wrap onto an extra line, so its row grows taller
$cup = $wb.Worksheets.Item("Cup")
$cup.Rows.Item(3).RowHeight = 76.5

# Cup ends up the active/visible tab
$cup.Activate()
